$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 108: only the timestamp in column A changes ---
$ws.Cells.Item(108, 1).Value = 45477.2916666667

# --- Row 109: brand new data row appended for BWZ.MI ---
$ws.Cells.Item(109, 1).Value = 45478.6353472222
# Reuse the exact same date/time style as the cell above (A108) so no new
# cell style gets introduced.
$ws.Cells.Item(108, 1).Copy()
$ws.Cells.Item(109, 1).PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(109, 2).Value = 28779
$ws.Cells.Item(109, 3).Value = 0.699999988079071
$ws.Cells.Item(109, 4).Value = 0.660000026226044
$ws.Cells.Item(109, 5).Value = 0.675000011920929
$ws.Cells.Item(109, 6).Value = 0.665000021457672

# Column G stores the adj_close figure as text (shared string), mirroring
# every other row in the sheet. Build it via a TEXT() formula so Excel
# treats the result as a genuine string instead of auto-converting the
# numeric-looking text back into a number, then collapse the formula down
# to its literal value (paste values only, keeping the default cell style).
$ws.Cells.Item(109, 7).Formula = '=TEXT(0.665000021457672,"0.000000000000000")'
$ws.Cells.Item(109, 7).Copy()
$ws.Cells.Item(109, 7).PasteSpecial(-4163)   # xlPasteValues

$ws.Cells.Item(109, 8).Value = "BWZ.MI"

$excel.CutCopyMode = $false
